# Refresh the cryptocurrency price/volume listing (Price = column D,
# Volume(1h) = column E) with the latest scraped values.
# A couple of D-column prices are plain numeric-looking strings
# (e.g. "1.001"); they are written with a leading apostrophe
# ('' inside a single-quoted PowerShell string yields one literal ')
# so Excel keeps them as text instead of silently parsing them as
# numbers, exactly like the original inline-string cells.
# Rows 20/21 additionally swap their Coin name / Link (B, C) because
# Avalanche and ShibaInu traded ranking positions in this refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '25.684.17'
$ws.Range("E2").Value = '  -2.93%  '
$ws.Range("D3").Value = '1.743.16'
$ws.Range("E3").Value = '  -5.03%  '
$ws.Range("D4").Value = '''1.001'
$ws.Range("E4").Value = '  -0.07%  '
$ws.Range("D5").Value = '''238.59'
$ws.Range("E5").Value = '  -8.36%  '
$ws.Range("E6").Value = '  -0.14%  '
$ws.Range("D7").Value = '''0.5044'
$ws.Range("E7").Value = '  -5.33%  '
$ws.Range("D8").Value = '''41.73'
$ws.Range("E8").Value = '  -6.99%  '
$ws.Range("D9").Value = '''0.2649'
$ws.Range("E9").Value = '  -11.93%  '
$ws.Range("D10").Value = '''0.06122'
$ws.Range("E10").Value = '  -10.75%  '
$ws.Range("D11").Value = '1.742.21'
$ws.Range("E11").Value = '  -5.16%  '
$ws.Range("D12").Value = '''0.06953'
$ws.Range("E12").Value = '  -4.73%  '
$ws.Range("D13").Value = '''15.22'
$ws.Range("E13").Value = '  -13.49%  '
$ws.Range("D14").Value = '''4.491'
$ws.Range("E14").Value = '  -9.46%  '
$ws.Range("D15").Value = '''0.5966'
$ws.Range("E15").Value = '  -18.88%  '
$ws.Range("D16").Value = '''76.53'
$ws.Range("E16").Value = '  -13.94%  '
$ws.Range("E17").Value = '  -0.13%  '
$ws.Range("E18").Value = '  -0.09%  '
$ws.Range("D19").Value = '25.700.02'
$ws.Range("E19").Value = '  -2.96%  '
$ws.Range("B20").Value = 'Avalanche'
$ws.Range("C20").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D20").Value = '''11.61'
$ws.Range("E20").Value = '  -16.45%  '
$ws.Range("B21").Value = 'ShibaInu'
$ws.Range("C21").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D21").Value = '''0.000006773'
$ws.Range("E21").Value = '  -14.22%  '
$ws.Range("D22").Value = '1.968.32'
$ws.Range("E22").Value = '  -5.39%  '
$ws.Range("D23").Value = '''4.039'
$ws.Range("E23").Value = '  -11.78%  '
$ws.Range("D24").Value = '''8.151'
$ws.Range("E24").Value = '  -11.76%  '
$ws.Range("D25").Value = '''5.107'
$ws.Range("E25").Value = '  -14.39%  '
$ws.Range("D26").Value = '''137.46'
$ws.Range("E26").Value = '  -3.75%  '
$ws.Range("D27").Value = '''1.517'
$ws.Range("E27").Value = '  -9.83%  '
$ws.Range("E28").Value = '  -17.40%  '
$ws.Range("E29").Value = '  -11.60%  '
$ws.Range("D30").Value = '''103.21'
$ws.Range("E30").Value = '  -6.41%  '
$ws.Range("D31").Value = '''3.755'
$ws.Range("E31").Value = '  -11.12%  '
$ws.Range("D32").Value = '''0.08100'
$ws.Range("E32").Value = '  -7.97%  '
$ws.Range("D33").Value = '''3.452'
$ws.Range("E33").Value = '  -13.99%  '
$ws.Range("D34").Value = '''0.04494'
$ws.Range("E34").Value = '  -6.21%  '
$ws.Range("D35").Value = '''0.9992'
$ws.Range("E35").Value = '  -0.18%  '
$ws.Range("D36").Value = '''2.652'
$ws.Range("E36").Value = '  -9.73%  '
$ws.Range("D37").Value = '''0.9816'
$ws.Range("E37").Value = '  -12.98%  '
$ws.Range("D38").Value = '''0.6092'
$ws.Range("E38").Value = '  -16.69%  '
$ws.Range("D39").Value = '''2.653'
$ws.Range("E39").Value = '  -14.22%  '
$ws.Range("D40").Value = '''0.01549'
$ws.Range("E40").Value = '  -9.24%  '
$ws.Range("D41").Value = '''1.913'
$ws.Range("E41").Value = '  -16.71%  '
$ws.Range("D42").Value = '''0.9998'
$ws.Range("E42").Value = '  -0.16%  '
$ws.Range("D43").Value = '''103.52'
$ws.Range("E43").Value = '  -3.64%  '
$ws.Range("D44").Value = '''0.3794'
$ws.Range("E44").Value = '  -19.51%  '
$ws.Range("D45").Value = '''5.116'
$ws.Range("E45").Value = '  -13.00%  '
$ws.Range("D46").Value = '''0.7264'
$ws.Range("E46").Value = '  -19.73%  '
$ws.Range("D47").Value = '''0.05333'
$ws.Range("E47").Value = '  -7.96%  '
$ws.Range("D48").Value = '''0.1111'
$ws.Range("E48").Value = '  -9.63%  '
$ws.Range("E49").Value = '  -13.49%  '
$ws.Range("D50").Value = '''5.882'
$ws.Range("E50").Value = '  -20.04%  '
$ws.Range("D51").Value = '''52.37'
$ws.Range("E51").Value = '  -12.85%  '
